$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 25228
$ws.Range("B2").Value = "Sra. Stella Viana"
$ws.Range("C2").Value = "Jurídico"
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 45100
$ws.Range("G2").Value = 6531.95

# Row 3
$ws.Range("A3").Value = 31023
$ws.Range("B3").Value = "Marcela da Cruz"
$ws.Range("C3").Value = "Financeiro"
$ws.Range("D3").Value = "Viagem de negócios"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 45081
$ws.Range("G3").Value = 6225.8

# Row 4
$ws.Range("A4").Value = 44280
$ws.Range("B4").Value = "Alexandre Cardoso"
$ws.Range("C4").Value = "Financeiro"
$ws.Range("D4").Value = "Problemas pessoais"
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 45091
$ws.Range("G4").Value = 11925.3

# Row 5
$ws.Range("A5").Value = 89716
$ws.Range("B5").Value = "Sra. Julia Martins"
$ws.Range("C5").Value = "Marketing"
$ws.Range("D5").Value = "Problemas pessoais"
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 45094
$ws.Range("G5").Value = 3704.8

# Row 6
$ws.Range("A6").Value = 23195
$ws.Range("B6").Value = "Nathan Nunes"
$ws.Range("C6").Value = "Marketing"
$ws.Range("D6").Value = "Outros"
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 45094
$ws.Range("G6").Value = 2832.78

# Row 7
$ws.Range("A7").Value = 50366
$ws.Range("B7").Value = "João Felipe Cavalcanti"
$ws.Range("C7").Value = "Atendimento ao Cliente"
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 45081
$ws.Range("G7").Value = 3791.59

# Row 8
$ws.Range("A8").Value = 1603
$ws.Range("B8").Value = "Clarice Pires"
$ws.Range("C8").Value = "Engenharia"
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 45099
$ws.Range("G8").Value = 6768.48

# Row 9
$ws.Range("A9").Value = 58846
$ws.Range("B9").Value = "Rafaela Martins"
$ws.Range("C9").Value = "Engenharia"
$ws.Range("D9").Value = "Consulta médica"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 45104
$ws.Range("G9").Value = 7483.06

# Row 10
$ws.Range("A10").Value = 24398
$ws.Range("B10").Value = "Isaac Duarte"
$ws.Range("C10").Value = "Engenharia"
$ws.Range("D10").Value = "Consulta médica"
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 45090
$ws.Range("G10").Value = 8072.1

# Row 11
$ws.Range("A11").Value = 31865
$ws.Range("B11").Value = "Cecília Araújo"
$ws.Range("C11").Value = "Atendimento ao Cliente"
$ws.Range("D11").Value = "Doença"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 45091
$ws.Range("G11").Value = 6539.18
